$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.722.02"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.601.86"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.81"
$ws.Range("D5").ClearFormats()
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "1.826.83"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "1.594.22"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.98"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "26.689.35"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "210.66"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.18"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.27"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.97"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.99"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").Value = "1.293.07"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.48"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.597"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.14"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.18%  "
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.830"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("D45").Value = "1.739.20"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.67"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0516"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.01"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.45"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.13%  "
